$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hours update (commit: "Updated hours for CJ and Nestor") ---

# Christopher Isherwood (CJ) - column D - Week 2 hours: 2 -> 6
$ws.Range("D5").Value = 6

# Nestor Macias - column F - hours logged for Week 2 and Week 3
$ws.Range("F5").Value = 12
$ws.Range("F6").Value = 2

# --- Update the active selection to match where the edits were made ---
$null = $ws.Range("G6").Select()
